# "Deal with the very first column being NULL"
#
# Adds a new row 6 to the "exotic" sheet that mirrors row 5's layout but
# leaves the first column (A) - and the already-always-empty D column -
# blank/NULL, renames the exotic underlined font away from the missing
# "Noto Sans CJK SC Regular" to "DejaVu Sans", and moves the live
# selection down to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Font used only by A5 ("Noto Sans CJK SC Regular", underlined) -> DejaVu Sans
$ws.Range("A5").Font.Name = "DejaVu Sans"

# New row 6: same shape as row 5, but column A (and D, already blank in
# every row) stay empty - "the very first column being NULL".
$ws.Range("B6").Value = -100
$ws.Range("C6").Formula = "=50+100-50"
$ws.Range("E6").Formula = "=FALSE()"
$ws.Range("E6").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("F6").Value = 54969
$ws.Range("F6").NumberFormat = "D"". ""MMM"". ""YYYY"

# Move the active selection to C10.
$ws.Range("C10").Select() | Out-Null
